# Update excess mortality scripts in line with rtweet update
#
# The upstream data-collection scripts were re-run (rtweet API update), which
# pulled slightly revised province-level weekly excess-mortality counts for
# weeks 10, 19, 24, 28, 29, 31-40 of 2022 and added the brand-new week 40 row.
# All percentage-change columns (AE:AP) are ROUND() formulas over the raw
# counts, so they recompute automatically once the raw counts are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: set a raw (non-formula) numeric cell -------------------------
function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 116 (2022 week 10)
Set-Cell "U116" 385

# Row 125 (2022 week 19)
Set-Cell "V125" 209

# Row 130 (2022 week 24)
Set-Cell "W130" 449

# Row 134 (2022 week 28)
Set-Cell "U134" 372

# Row 135 (2022 week 29)
Set-Cell "X135" 639

# Row 137 (2022 week 31)
Set-Cell "W137" 455
Set-Cell "X137" 614
Set-Cell "Z137" 446

# Row 138 (2022 week 32)
Set-Cell "X138" 596

# Row 139 (2022 week 33)
Set-Cell "X139" 604

# Row 140 (2022 week 34)
Set-Cell "W140" 443

# Row 141 (2022 week 35)
Set-Cell "W141" 412
Set-Cell "X141" 597
Set-Cell "Z141" 477

# Row 142 (2022 week 36)
Set-Cell "U142" 370
Set-Cell "W142" 454

# Row 143 (2022 week 37)
Set-Cell "Q143" 110
Set-Cell "S143" 193
Set-Cell "U143" 374
Set-Cell "W143" 405
Set-Cell "X143" 601
Set-Cell "AA143" 223

# Row 144 (2022 week 38)
Set-Cell "S144" 205
Set-Cell "U144" 376
Set-Cell "V144" 227
Set-Cell "W144" 479
Set-Cell "X144" 602
Set-Cell "Z144" 424
Set-Cell "AA144" 218

# Row 145 (2022 week 39) - revised set of 2022 counts
Set-Cell "P145" 120
Set-Cell "Q145" 122
Set-Cell "R145" 119
Set-Cell "S145" 223
Set-Cell "T145" 34
Set-Cell "U145" 435
Set-Cell "V145" 231
Set-Cell "W145" 472
Set-Cell "X145" 648
Set-Cell "Y145" 75
Set-Cell "Z145" 429
# AA145 unchanged (248)

# Row 146 (2022 week 40) - newly available week; fill in the 2022 figures and
# the corresponding "change vs. baseline" formulas (columns N:AP), mirroring
# the pattern already used by the preceding rows.
Set-Cell "N146" 2022
Set-Cell "O146" 40
Set-Cell "P146" 134
Set-Cell "Q146" 131
Set-Cell "R146" 104
Set-Cell "S146" 247
Set-Cell "T146" 47
Set-Cell "U146" 454
Set-Cell "V146" 198
Set-Cell "W146" 456
Set-Cell "X146" 611
Set-Cell "Y146" 81
Set-Cell "Z146" 474
Set-Cell "AA146" 226
Set-Cell "AC146" 2022
Set-Cell "AD146" 40

$ws.Range("AE146").Formula = "=ROUND((P146-B146)/B146*100,2)"
$ws.Range("AF146").Formula = "=ROUND((Q146-C146)/C146*100,2)"
$ws.Range("AG146").Formula = "=ROUND((R146-D146)/D146*100,2)"
$ws.Range("AH146").Formula = "=ROUND((S146-E146)/E146*100,2)"
$ws.Range("AI146").Formula = "=ROUND((T146-F146)/F146*100,2)"
$ws.Range("AJ146").Formula = "=ROUND((U146-G146)/G146*100,2)"
$ws.Range("AK146").Formula = "=ROUND((V146-H146)/H146*100,2)"
$ws.Range("AL146").Formula = "=ROUND((W146-I146)/I146*100,2)"
$ws.Range("AM146").Formula = "=ROUND((X146-J146)/J146*100,2)"
$ws.Range("AN146").Formula = "=ROUND((Y146-K146)/K146*100,2)"
$ws.Range("AO146").Formula = "=ROUND((Z146-L146)/L146*100,2)"
$ws.Range("AP146").Formula = "=ROUND((AA146-M146)/M146*100,2)"

# Recalculate all formulas (AE:AP ROUND() columns) so cached <v> values in the
# saved workbook reflect the corrected raw counts above.
$excel.CalculateFullRebuild()

# Restore the last on-screen selection recorded for this sheet.
$ws.Range("AE138").Select()
